# feat: add 2022-Q3 data
#
# 1. Insert a new worksheet "2022-Q3" right after "总计" (before "2022-Q2")
#    and populate it with the Q3 fund-holding detail table.
# 2. Insert a new row at the top of the "总计" summary table's data
#    (row 2) for the 2022-Q3 totals, pushing the existing quarters down.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) New sheet "2022-Q3" inserted after "总计" (position 2)
# ---------------------------------------------------------------------
# Duplicate the "2022-Q2" sheet (same column layout/header/formatting)
# right before itself, then rename + overwrite its data — this keeps the
# header styling / borders / column widths consistent with the other
# quarterly sheets without having to rebuild formatting from scratch.
$template = $wb.Worksheets.Item(2)
$template.Copy($wb.Worksheets.Item(2))
$qsheet = $wb.Worksheets.Item(2)
$qsheet.Name = "2022-Q3"

# The template ("2022-Q2") has 8 data rows; the Q3 table only needs 6, so
# drop the last two leftover rows.
$qsheet.Rows.Item(8).Delete()
$qsheet.Rows.Item(7).Delete()

# Data rows: index, code, name, scale, position%, ratio%, market value, rank
$qdata = @(
    @(0, "002560", "诺安和鑫灵活配置混合",       "32.70", "79.56", "7.00", "2.2890", 6),
    @(1, "320022", "诺安研究精选股票",           "6.17",  "92.67", "2.16", "0.1333", 8),
    @(2, "001706", "诺安积极回报灵活配置混合A",  "0.88",  "94.16", "9.55", "0.0840", 2),
    @(3, "012847", "诺安积极回报灵活配置混合C",  "0.40",  "94.16", "9.55", "0.0382", 2),
    @(4, "007533", "格林创新成长混合A",          "0.05",  "88.63", "4.53", "0.0023", 10),
    @(5, "007534", "格林创新成长混合C",          "0.04",  "88.63", "4.53", "0.0018", 10)
)

for ($i = 0; $i -lt $qdata.Count; $i++) {
    $r = $i + 2
    $row = $qdata[$i]
    $qsheet.Cells.Item($r, 1).Value = $row[0]
    $qsheet.Cells.Item($r, 2).Value = "'" + $row[1]
    $qsheet.Cells.Item($r, 3).Value = $row[2]
    $qsheet.Cells.Item($r, 4).Value = "'" + $row[3]
    $qsheet.Cells.Item($r, 5).Value = "'" + $row[4]
    $qsheet.Cells.Item($r, 6).Value = "'" + $row[5]
    $qsheet.Cells.Item($r, 7).Value = "'" + $row[6]
    $qsheet.Cells.Item($r, 8).Value = $row[7]
}

# ---------------------------------------------------------------------
# 2) "总计" (totals) sheet: insert a new data row for 2022-Q3 at row 2
# ---------------------------------------------------------------------
$total = $wb.Worksheets.Item(1)
$total.Rows.Item(2).Insert()

# Row 2 picks up row-1 (header) formatting on insert; reset it to match
# the plain data-row formatting used by the rest of the table (copy from
# what is now row 3, the old row 2).
$total.Range("A3:D3").Copy()
$total.Range("A2:D2").PasteSpecial(-4122)

$total.Cells.Item(2, 1).Value = 0
$total.Cells.Item(2, 2).Value = "2022-Q3"
$total.Cells.Item(2, 3).Value = 6
$total.Cells.Item(2, 4).Value = 2.55

# Renumber the index column for the rows that shifted down.
$total.Cells.Item(3, 1).Value = 1
$total.Cells.Item(4, 1).Value = 2
$total.Cells.Item(5, 1).Value = 3
$total.Cells.Item(6, 1).Value = 4
